# ---------------------------------------------------------------------------
# Applies the "BOLETA HOY" batch of new invoice rows to
# "FORMATO PARA LLENADO DE FACTURAS Y BOLETAS - VILLAFUERTE - 5 ESQUINAS.xlsx"
#
# Net effect (see commit diff):
#   * Rows 145/146 (previously a BOLETA/EB01 placeholder line) become blank
#     spacer rows.
#   * Rows 147-161 get the 15 new transactions that were typed in (dated
#     2024-01-25 / serial 45316), replacing what used to be empty template
#     rows.
#   * Rows 162/163 become blank spacer rows (gain a formatted, empty G cell).
#   * Rows 164-167 become four more BOLETA/EB01 placeholder rows (G formula
#     evaluates to 0 because E/F are blank), replacing old placeholder rows
#     that contained stray 45313 dates.
#   * Two new blank rows (2053/2054) are appended at the very bottom so the
#     sheet's used range grows from A3:J2052 to A3:J2054.
#   * The shared formula in G84:G... and G140:G... is re-anchored because
#     the formula now spans further down (through row 167).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$FORMULA_TEMPLATE = "=IF(ISBLANK(E{0}),F{0},E{0}*1.18)"

# --- Rows 145 & 146: clear back down to blank spacer rows -----------------
# They used to hold a BOLETA/EB01 row with nothing typed in (B/C only) --
# remove those two cells outright and blank out A/G while keeping their
# existing number formats (date / "0.0").
$ws.Range("B145:C146").Clear()
$ws.Range("A145:A146").ClearContents()
$ws.Range("G145:G146").ClearContents()

# --- Rows 147-161: the 15 new transactions ---------------------------------
# columns: Date(A) Doc(B) Serie(C) Numero(D) Gravado(E) Exonerado(F) Total(G) Cliente(H)
$rows = @(
    @{ R=147; D=14033; E=917.8;   F=$null; H="ERWIN ROBERT YAYA MANCO" },
    @{ R=148; D=14035; E=$null;   F=701;   H="KATHERINE TAMARIZ BERNAL" },
    @{ R=149; D=14042; E=$null;   F=1260;  H="SAMUEL MAXIMO QUISPE CONDEZO" },
    @{ R=150; D=14043; E=227.97;  F=$null; H="SOLEDAD MILAGROS PADIN RAMOS" },
    @{ R=151; D=14047; E=$null;   F=700;   H="JORGE LUIS SURICHAQUI ACOSTA" },
    @{ R=152; D=14050; E=610.17;  F=$null; H="JUSTA FERMINA LAURA DE CHUQUISPUMA" },
    @{ R=153; D=14052; E=1609.32; F=$null; H="RUBEN FLORIAN ZAPATA" },
    @{ R=155; D=14053; E=537.29;  F=$null; H="CATALINO CORONADO OLIVERA" },
    @{ R=156; D=14054; E=33.9;    F=$null; H="GILBERTO AUGUSTO PAREJA TORRES" },
    @{ R=157; D=14055; E=$null;   F=700;   H="MIGUEL RODRIGUEZ QUISPE" },
    @{ R=158; D=14056; E=$null;   F=700;   H="PILAR ESTHER FLORIAN ZAPATA" },
    @{ R=159; D=$null; E=$null;   F=700;   H="ROMULO GUTIERREZ QUISPE" },
    @{ R=160; D=$null; E=$null;   F=700;   H="RUBEN FLORIAN ZAPATA" },
    @{ R=161; D=$null; E=$null;   F=700;   H="JOSE LUIS LUYO SANCHEZ" }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = 45316
    $ws.Cells.Item($r, 2).Value = "BOLETA"
    $ws.Cells.Item($r, 3).Value = "EB01"
    if ($row.D -ne $null) {
        $ws.Cells.Item($r, 4).Value = $row.D
    }
    if ($row.E -ne $null) {
        $ws.Cells.Item($r, 5).Value = $row.E
    }
    if ($row.F -ne $null) {
        $ws.Cells.Item($r, 6).Value = $row.F
    }
    $ws.Cells.Item($r, 7).NumberFormat = "0.0"
    $ws.Cells.Item($r, 7).Formula = [string]::Format($FORMULA_TEMPLATE, $r)
    $ws.Cells.Item($r, 8).Value = $row.H
}

# Row 154 is a FACTURA/E001 line (AGRO RUNA) rather than a BOLETA/EB01 one.
$r = 154
$ws.Cells.Item($r, 1).Value = 45316
$ws.Cells.Item($r, 2).Value = "FACTURA"
$ws.Cells.Item($r, 3).Value = "E001"
$ws.Cells.Item($r, 4).Value = 3157
$ws.Cells.Item($r, 5).Value = 61.02
$ws.Cells.Item($r, 7).NumberFormat = "0.0"
$ws.Cells.Item($r, 7).Formula = [string]::Format($FORMULA_TEMPLATE, $r)
$ws.Cells.Item($r, 8).Value = "AGRO RUNA S.A.C. "

# Row 160's client cell keeps the odd one-off Calibri-10pt look that an
# existing cell elsewhere in the sheet (H119) already uses -- copy that
# format across instead of re-describing the font, so we reuse the same
# style entry rather than minting a near-duplicate font.
$ws.Range("H119").Copy()
$ws.Cells.Item(160, 8).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Rows 162 & 163: new blank spacer rows (gain a formatted G cell) ------
$ws.Range("G162:G163").NumberFormat = "0.0"

# --- Rows 164-167: four more BOLETA/EB01 placeholder rows -----------------
foreach ($r in 164..167) {
    $ws.Cells.Item($r, 1).Value = 45316
    $ws.Cells.Item($r, 2).Value = "BOLETA"
    $ws.Cells.Item($r, 3).Value = "EB01"
    $ws.Cells.Item($r, 7).NumberFormat = "0.0"
    $ws.Cells.Item($r, 7).Formula = [string]::Format($FORMULA_TEMPLATE, $r)
}

# --- Two new blank rows appended at the very end of the sheet --------------
# Copy A2052's (empty) format down so the new rows pick up the same style
# entry instead of minting a new one via NumberFormat.
$ws.Range("A2052").Copy()
$ws.Range("A2053").PasteSpecial(-4122)
$ws.Range("A2054").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(2053).RowHeight = $ws.Rows.Item(2052).RowHeight
$ws.Rows.Item(2054).RowHeight = $ws.Rows.Item(2052).RowHeight

# --- View state: keep the frozen pane, move the live selection -----------
$ws.Range("F163").Select()
